# edit.ps1 - applies the changes described by the target diff:
#  1) Title block: merge the multi-run "ADNAN HAIDER"/" - SP20-BSE-0"/"37"
#     paragraph into a single run "ADNAN HAIDER - SP20-BSE-037", and merge
#     the two-run "LOG IN"/" (USE CASE)" paragraph into a single run
#     "LOG IN (USE CASE)".
#  2) Post-conditions paragraph: replace the
#     'Login instance "LIA" has been created ... LIA was associated ...'
#     wording with the generic 'Login instance has been created ...
#     Login instance was associated ...' wording, and drop a _GoBack
#     bookmark right before "was associated with the User" (this mirrors
#     where Word leaves the _GoBack bookmark after the last edit).

$d = $word.ActiveDocument

# Useful characters that are awkward to embed literally.
$enDash = [char]0x2013
$ldquo  = [char]0x201C
$rdquo  = [char]0x201D

# --- 1a. "ADNAN HAIDER" + " - SP20-BSE-0" + "37"  ->  "ADNAN HAIDER - SP20-BSE-037"
$oldTitle1 = "ADNAN HAIDER" + " " + $enDash + " SP20-BSE-0" + "37"
$newTitle1 = "ADNAN HAIDER" + " " + $enDash + " SP20-BSE-037"
$d.Content.Find.Execute($oldTitle1, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newTitle1, 2) | Out-Null

# --- 1b. "LOG IN" + " (USE CASE)"  ->  "LOG IN (USE CASE)" (single run)
$oldTitle2 = "LOG IN (USE CASE)"
$newTitle2 = "LOG IN (USE CASE)"
$d.Content.Find.Execute($oldTitle2, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newTitle2, 2) | Out-Null

# --- 2. Post conditions wording
$oldPost = "Login instance " + $ldquo + "LIA" + $rdquo + `
           " has been created for User login operation. LIA was associated with the User "
$newPost = "Login instance has been created for User login operation. " + `
           "Login instance was associated with the User "
$d.Content.Find.Execute($oldPost, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newPost, 2) | Out-Null

# Drop a _GoBack bookmark right before "was associated with the User" - this
# is where Word records the last editing position.
$bmRange = $d.Content
$bmRange.Find.Execute("was associated with the User", $true, $false, $false, `
                       $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
